$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Heures TD" column (F) to hold "Groupes CM".
# This shifts Heures TD / Groupes TD / Heures TP / Groupes TP one column to the right.
$ws.Columns.Item(6).Insert()

# New column F header and values ("Groupes CM" = 1 group of CM for every course)
$ws.Cells.Item(1, 6).Value = "Groupes CM"
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Refresh the worksheet's stored sort state so its range keeps up with the
# newly added column (A2:I33 -> A2:J33), preserving the original sort keys.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B2:B33")) | Out-Null
$sort.SortFields.Add($ws.Range("A2:A33")) | Out-Null
$sort.SetRange($ws.Range("A2:J33"))
$sort.Header = 2
$sort.Apply()

# Leave the selection where the edit ended up.
$ws.Range("D25").Select()
